# Dev IV Project Rubric - retake edit
# Student switched their "Student(I, II, or III)" milestone marker for every
# completed feature row to "I" (this retake only covers Milestone I), removed
# the "Milestone II Complete" and "Milestone III Complete" markers for the
# "Effective Use of GIT" row, and updated the Git repo address + selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- A2: update the student's Git repository address ---
$ws.Range("A2").Value = "Student Git Address:  https://github.com/maximusrex17/Project4"

# --- Column E: which milestone (I, II, or III) each completed feature belongs to ---
# Rows that changed milestone to "I":
$ws.Range("E4").Value = "I"
$ws.Range("E5").Value = "I"
$ws.Range("E8").Value = "I"
$ws.Range("E9").Value = "I"
$ws.Range("E15").Value = "I"
$ws.Range("E24").Value = "I"
$ws.Range("E57").Value = "I"
$ws.Range("E63").Value = "I"

# Row 6 feature is no longer marked complete for any milestone - clear it.
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

# --- Effective Use of GIT (row 91) and cleanup row (row 92): only Milestone I now ---
$ws.Range("D91").ClearContents()
$ws.Range("E91").ClearContents()
$ws.Range("D92").ClearContents()
$ws.Range("E92").ClearContents()

# --- Update view: scroll down and move the active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E6").Select()
